# Applies the stock-quantity / value corrections described in the commit diff.
# For every affected item row, column F = quantity and column G = value (|= D * F|).
# A handful of rows were re-sequenced (their whole row of data swapped with the
# row immediately below/above them) - those are written as two explicit row updates.
# "Sub Total:" / "Grand Total:" rows in column B are plain numbers (not formulas in
# this workbook), so each subtotal affected by the item-level edits is corrected too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F44").Value = 404
$ws.Range("G44").Value = 14709.64

$ws.Range("F47").Value = 190
$ws.Range("G47").Value = 36649.1

$ws.Range("F54").Value = 224
$ws.Range("G54").Value = 12566.4

$ws.Range("F58").Value = 11
$ws.Range("G58").Value = 649

$ws.Range("F60").Value = 47
$ws.Range("G60").Value = 2096.2

$ws.Range("F61").Value = 108
$ws.Range("G61").Value = 6022.08

$ws.Range("F64").Value = 49
$ws.Range("G64").Value = 3818.57

$ws.Range("B72").Value = 178619.98

$ws.Range("F120").Value = 17
$ws.Range("G120").Value = 794.58

$ws.Range("F126").Value = 80
$ws.Range("G126").Value = 10781.6

$ws.Range("B129").Value = 68434.10000000001

# Row 132 (re-sequenced item)
$ws.Range("B132").Value = 65258
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0

# Row 133 (re-sequenced item)
$ws.Range("B133").Value = 64196
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 32143.58

$ws.Range("F145").Value = 27
$ws.Range("G145").Value = 1839.51

$ws.Range("F147").Value = 36
$ws.Range("G147").Value = 1874.16

$ws.Range("F152").Value = 12
$ws.Range("G152").Value = 788.04

$ws.Range("B153").Value = 19409.03

$ws.Range("F178").Value = 87
$ws.Range("G178").Value = 5526.24

$ws.Range("F179").Value = 33
$ws.Range("G179").Value = 2594.13

$ws.Range("F188").Value = 10
$ws.Range("G188").Value = 895.4

$ws.Range("B199").Value = 57190.15

$ws.Range("F216").Value = 97
$ws.Range("G216").Value = 5277.77

$ws.Range("F218").Value = 48
$ws.Range("G218").Value = 4276.32

$ws.Range("F221").Value = 149
$ws.Range("G221").Value = 16731.21

$ws.Range("B224").Value = 66971.11

$ws.Range("F228").Value = 370
$ws.Range("G228").Value = 6845

$ws.Range("B235").Value = 15532.59

$ws.Range("F261").Value = 5
$ws.Range("G261").Value = 1579

$ws.Range("F262").Value = 32
$ws.Range("G262").Value = 2580.48

$ws.Range("F265").Value = 44
$ws.Range("G265").Value = 4617.36

$ws.Range("F269").Value = 3
$ws.Range("G269").Value = 311.1

$ws.Range("F278").Value = 17
$ws.Range("G278").Value = 4626.04

$ws.Range("F283").Value = 2
$ws.Range("G283").Value = 100.82

$ws.Range("F284").Value = 27
$ws.Range("G284").Value = 3659.58

$ws.Range("F286").Value = 10
$ws.Range("G286").Value = 971

$ws.Range("F291").Value = 2
$ws.Range("G291").Value = 221.88

$ws.Range("F297").Value = 20
$ws.Range("G297").Value = 1713.6

# Row 298 (re-sequenced item)
$ws.Range("B298").Value = 66196
$ws.Range("C298").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F298").Value = 5
$ws.Range("G298").Value = 438.5

# Row 299 (re-sequenced item)
$ws.Range("B299").Value = 64985
$ws.Range("C299").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F299").Value = 12
$ws.Range("G299").Value = 1052.4

$ws.Range("B301").Value = 100263.71

# Row 310 (re-sequenced item)
$ws.Range("B310").Value = 55373
$ws.Range("E310").Value = 163.62
$ws.Range("F310").Value = -94
$ws.Range("G310").Value = -13562.32

# Row 311 (re-sequenced item)
$ws.Range("B311").Value = 63520
$ws.Range("E311").Value = 153.4
$ws.Range("F311").Value = 35
$ws.Range("G311").Value = 5049.8

$ws.Range("B334").Value = -22661.27

$ws.Range("F351").Value = 47
$ws.Range("G351").Value = 7712.7

$ws.Range("F355").Value = 131
$ws.Range("G355").Value = 9776.530000000001

$ws.Range("B362").Value = 74265.56

$ws.Range("F367").Value = 202
$ws.Range("G367").Value = 28399.18

$ws.Range("B369").Value = 64210.3

$ws.Range("F374").Value = 45
$ws.Range("G374").Value = 1440.9

$ws.Range("F376").Value = 174
$ws.Range("G376").Value = 28882.26

$ws.Range("B378").Value = 50505.68

$ws.Range("F391").Value = 30
$ws.Range("G391").Value = 892.2

$ws.Range("F393").Value = 370
$ws.Range("G393").Value = 35742

$ws.Range("B395").Value = 51832.06

$ws.Range("F402").Value = 116
$ws.Range("G402").Value = 2955.68

$ws.Range("F403").Value = 63
$ws.Range("G403").Value = 2278.71

$ws.Range("F408").Value = 19
$ws.Range("G408").Value = 651.89

$ws.Range("F409").Value = 63
$ws.Range("G409").Value = 2554.02

$ws.Range("F414").Value = 170
$ws.Range("G414").Value = 2694.5

$ws.Range("F419").Value = 68
$ws.Range("G419").Value = 3914.76

$ws.Range("F422").Value = 50
$ws.Range("G422").Value = 1468

$ws.Range("B423").Value = 158363.79

$ws.Range("F437").Value = 7
$ws.Range("G437").Value = 188.23

$ws.Range("B444").Value = 21802.16

$ws.Range("F460").Value = 57
$ws.Range("G460").Value = 16128.72

$ws.Range("B464").Value = 84773.05

# Row 502 (re-sequenced item)
$ws.Range("B502").Value = 64833
$ws.Range("E502").Value = 34.9
$ws.Range("F502").Value = 88
$ws.Range("G502").Value = 2889.04

# Row 503 (re-sequenced item)
$ws.Range("B503").Value = 60025
$ws.Range("E503").Value = 37.22
$ws.Range("F503").Value = -98
$ws.Range("G503").Value = -3217.34

# Row 512 (re-sequenced item)
$ws.Range("B512").Value = 60022
$ws.Range("E512").Value = 37.22
$ws.Range("F512").Value = -113
$ws.Range("G512").Value = -3709.79

# Row 513 (re-sequenced item)
$ws.Range("B513").Value = 64830
$ws.Range("E513").Value = 34.9
$ws.Range("F513").Value = 83
$ws.Range("G513").Value = 2724.89

$ws.Range("F517").Value = 185
$ws.Range("G517").Value = 18475.95

$ws.Range("F518").Value = 12
$ws.Range("G518").Value = 1422.96

$ws.Range("F527").Value = 75
$ws.Range("G527").Value = 2055

$ws.Range("F528").Value = 42
$ws.Range("G528").Value = 1118.88

$ws.Range("B531").Value = 109906.03

$ws.Range("F533").Value = 26
$ws.Range("G533").Value = 860.86

$ws.Range("F535").Value = 107
$ws.Range("G535").Value = 3542.77

$ws.Range("F537").Value = 182
$ws.Range("G537").Value = 6026.02

$ws.Range("B541").Value = 20017.79

$ws.Range("F564").Value = 146
$ws.Range("G564").Value = 17790.1

$ws.Range("B567").Value = 20129.44

$ws.Range("F611").Value = 161
$ws.Range("G611").Value = 21429.1

$ws.Range("B613").Value = 21429.1

$ws.Range("F616").Value = 56
$ws.Range("G616").Value = 1418.48

$ws.Range("F618").Value = 221
$ws.Range("G618").Value = 33240.61

$ws.Range("F631").Value = 296
$ws.Range("G631").Value = 10901.68

$ws.Range("B634").Value = 195376.77

$ws.Range("F669").Value = 35
$ws.Range("G669").Value = 1251.95

$ws.Range("B674").Value = 9991.49

$ws.Range("F680").Value = 475
$ws.Range("G680").Value = 77477.25

$ws.Range("B686").Value = 78489.8

$ws.Range("F702").Value = 3
$ws.Range("G702").Value = 1310.1

$ws.Range("B719").Value = 58929.68

$ws.Range("B724").Value = 2411053.84

$ws.Range("B725").Value = 2411053.84
